$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Strip the "root-code-name." prefix from the code-name cells in column A.
$ws.Range("A3").Value  = "common-code-name-1"
$ws.Range("A4").Value  = "common-code-name-1.common-code-name-2"
$ws.Range("A5").Value  = "common-code-name-1.common-code-name-2.common-code-name-3-1"
$ws.Range("A6").Value  = "common-code-name-1.common-code-name-2.common-code-name-3-2"
$ws.Range("A7").Value  = "common-code-name-1.common-code-name-2.common-code-name-3-3"
$ws.Range("A8").Value  = "error-code-name-1"
$ws.Range("A9").Value  = "error-code-name-1.error-code-name-2"
$ws.Range("A10").Value = "error-code-name-1.error-code-name-2.error-code-name-3-1"
$ws.Range("A11").Value = "error-code-name-1.error-code-name-2.error-code-name-3-2"
$ws.Range("A12").Value = "error-code-name-1.error-code-name-2.error-code-name-3-3"
$ws.Range("A13").Value = "error-2-code-name-1"
$ws.Range("A14").Value = "error-2-code-name-1.error-2-code-name-2"
$ws.Range("A15").Value = "error-2-code-name-1.error-2-code-name-2.error-2-code-name-3-1"
$ws.Range("A16").Value = "error-2-code-name-1.error-2-code-name-2.error-2-code-name-3-2"
$ws.Range("A17").Value = "error-2-code-name-1.error-2-code-name-2.error-2-code-name-3-3"

# Move the active cell/selection from B18 to A4.
$ws.Range("A4").Select()
